# Apply habitat_confirmations_priorities changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: Unnamed Tributary to Elk River / 197555_us -> hab_value = high, upstream_habitat_length_m = 6000
$ws.Range("H2").Value = "high"
$ws.Range("M2").Value = 6000

# Row 3: Unnamed Tributary to Elk River / 197555_ds -> hab_value = high
$ws.Range("H3").Value = "high"

# Row 8: Brule Creek / 197533_us -> upstream_habitat_length_m corrected from 0.125 to 125
$ws.Range("M8").Value = 125

# Row 15: Unnamed Tributary to Lizard Creek / 50159_us -> hab_value = high, upstream_habitat_length_m = 350
$ws.Range("H15").Value = "high"
$ws.Range("M15").Value = 350

# Row 16: Unnamed Tributary to Lizard Creek / 50159_ds -> hab_value = high
$ws.Range("H16").Value = "high"

# Row 17: Unnamed Tributary to Lizard Creek / 50155_us -> upstream_habitat_length_m = 1800
$ws.Range("M17").Value = 1800

# Row 21: Unnamed Tributary to Morrissey Creek / 50185_us -> hab_value = high, upstream_habitat_length_m = 4500
$ws.Range("H21").Value = "high"
$ws.Range("M21").Value = 4500

# Row 22: Unnamed Tributary to Morrissey Creek / 50185_ds -> hab_value = high
$ws.Range("H22").Value = "high"

# Update the selection / view state to match the edited workbook
$ws.Activate()
$ws.Range("I12").Select()
$excel.ActiveWindow.ScrollColumn = 2
